$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 7, shifting rows 7-12 down to 8-13.
$ws.Range("A7:T7").Insert()

# Copy formatting/style of old row 7 (now row 8, column A-T only) into the
# newly inserted row 7, so the date style (s="2" on column D) carries over.
$ws.Range("A8:T8").Copy()
$ws.Range("A7:T7").PasteSpecial(-4122) | Out-Null

# Populate the new row 7 with its data (same template values as the row below,
# with updated D, N, O, P, Q, S columns per the diff).
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 45030
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100107
$ws.Cells.Item(7, 8).Value = "Otros"
$ws.Cells.Item(7, 9).Value = 100107011
$ws.Cells.Item(7, 10).Value = "Tuna"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 16000
$ws.Cells.Item(7, 16).Value = 15500
$ws.Cells.Item(7, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(7, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(7, 19).Value = 861
$ws.Cells.Item(7, 20).Value = 18
